$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.053021666666666
$ws.Range("H2").Value = 21.159065
$ws.Range("I2").Value = 0.08011112358180576
$ws.Range("J2").Value = 0.08011112358180575
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.03637433333333333
$ws.Range("N2").Value = 0.109123
$ws.Range("Q2").Value = 0.2565489611105555
$ws.Range("R2").Value = 2.308940649995
$ws.Range("S2").Value = 0.08011112358180576
$ws.Range("T2").Value = 0.08011112358180575

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 35.63076933333333
$ws.Range("H3").Value = 106.892308
$ws.Range("I3").Value = 0.4047089460773642
$ws.Range("J3").Value = 0.4047089460773642
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.03637433333333333
$ws.Range("N3").Value = 0.109123
$ws.Range("Q3").Value = 1.296045480653778
$ws.Range("R3").Value = 11.664409325884
$ws.Range("S3").Value = 0.4047089460773642
$ws.Range("T3").Value = 0.4047089460773642

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 33.64714066666667
$ws.Range("H4").Value = 100.941422
$ws.Range("I4").Value = 0.3821780750881576
$ws.Range("J4").Value = 0.3821780750881575
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.03637433333333333
$ws.Range("N4").Value = 0.109123
$ws.Range("Q4").Value = 1.223892310322889
$ws.Range("R4").Value = 11.015030792906
$ws.Range("S4").Value = 0.3821780750881576
$ws.Range("T4").Value = 0.3821780750881575

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.048222
$ws.Range("H5").Value = 24.144666
$ws.Range("I5").Value = 0.09141501865831142
$ws.Range("J5").Value = 0.09141501865831139
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.03637433333333333
$ws.Range("N5").Value = 0.109123
$ws.Range("Q5").Value = 0.2927487097686667
$ws.Range("R5").Value = 2.634738387918
$ws.Range("S5").Value = 0.09141501865831142
$ws.Range("T5").Value = 0.09141501865831139

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.661325
$ws.Range("H6").Value = 10.983975
$ws.Range("I6").Value = 0.04158683659436109
$ws.Range("J6").Value = 0.04158683659436108
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.03637433333333333
$ws.Range("N6").Value = 0.109123
$ws.Range("Q6").Value = 0.1331782559916667
$ws.Range("R6").Value = 1.198604303925
$ws.Range("S6").Value = 0.04158683659436109
$ws.Range("T6").Value = 0.04158683659436108

